# Update Regresi API Testing 16/07
# Updates the error-message reference data on Sheet1:
#  - F9 "Register Failed" -> full dialog text (message + button), wrap text
#  - D10 email updated from sal3@gmail.com to sal5@gmail.com
#  - F10 "Register Success" -> full dialog text (message + button), wrap text
#  - Row heights for rows 9/10 grow to fit the wrapped text
#  - Selection moved to G3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# F9: Register Failed dialog text, now wrapped across 3 lines.
$ws.Range("F9").Value = "Register Failed" + $nl + "Account already exist , please try using another account" + $nl + "OK"
$ws.Range("F9").WrapText = $true

# D10: test account email changed.
$ws.Range("D10").Value = "sal5@gmail.com"

# F10: Register Success dialog text, now wrapped across 3 lines.
$ws.Range("F10").Value = "Register Success" + $nl + "Please login using your new account" + $nl + "Login"
$ws.Range("F10").WrapText = $true

# Rows grow to fit the new 3-line wrapped text.
$ws.Range("9:9").RowHeight = 45
$ws.Range("10:10").RowHeight = 45

# Active selection moved to G3.
$null = $ws.Range("G3").Select()
